$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# The single paragraph "Iterate through all variables against one or two
# (attrition/jobSatisfaction)" (originally split across two runs) is
# replaced with "Figure out which geoms to use with different types of
# data", and a brand-new sibling bullet paragraph (same list formatting)
# is inserted right after it holding the original sentence as one run.

$oldText = "Iterate through all variables against one or two (attrition/jobSatisfaction)"
$newText = "Figure out which geoms to use with different types of data"

$rng = $d.Content
$rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newText, 2) | Out-Null

# InsertParagraphAfter adds a fresh empty paragraph right after $rng,
# copying the paragraph/run formatting ($rng now points at the replaced
# sentence, so the new paragraph inherits that bullet's pPr/rPr).
$rng.InsertParagraphAfter() | Out-Null

# Paragraph.Range.Text (unlike a plain Range) includes the trailing
# paragraph mark, so trim it off before comparing.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($t -eq $newText) {
        $targetIndex = $i
        break
    }
}
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = $oldText

# --- Change 2 -------------------------------------------------------------
# "*Mon " + "5/18 - Wed 5/20 (afternoon)" (two runs) collapse into a
# single run "*Mon 5/18 - Wed 5/20 (afternoon)".

$monText = "*Mon 5/18 - Wed 5/20 (afternoon)"
$d.Content.Find.Execute($monText, $true, $false, $false, $false, $false, `
    $true, 1, $false, $monText, 2) | Out-Null
